$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Name of the Company(as per certificate)*"
$ws.Range("D1").Value = "Current Name of the Company*"
$ws.Range("E1").Value = "ISIN*"
$ws.Range("G1").Value = "Face Value*"
$ws.Range("H1").Value = "Closing Price in NSE*"
$ws.Range("I1").Value = "Closing Price in BSE*"
